$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.629.59"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.640.88"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "0.0625"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "19.08"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.867.65"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.666.28"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "4.17"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "64.85"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "26.643.27"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "215.12"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "9.45"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +12.41%  "
$ws.Range("D25").Value = "144.99"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "0.0512"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "1.277.06"
$ws.Range("E34").Value = "  +5.38%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").Value = "0.530"
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("D39").Value = "0.824"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "0.809"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D44").Value = "1.778.32"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "91.11"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").Value = "59.18"
$ws.Range("E46").Value = "  +8.35%  "
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").Value = "0.406"
$ws.Range("E51").Value = "  -0.53%  "
